$d = $word.ActiveDocument

$d.Content.Find.Execute("495÷4=123, 3", $true, $false, $false, $false, $false, $true, 1, $false, "585÷8=73, 1", 2) | Out-Null
$d.Content.Find.Execute("735÷4=183, 3", $true, $false, $false, $false, $false, $true, 1, $false, "978÷9=108, 6", 2) | Out-Null
$d.Content.Find.Execute("572÷3=190, 2", $true, $false, $false, $false, $false, $true, 1, $false, "655÷4=163, 3", 2) | Out-Null
$d.Content.Find.Execute("162÷8=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "532÷2=266, 0", 2) | Out-Null
$d.Content.Find.Execute("747÷3=249, 0", $true, $false, $false, $false, $false, $true, 1, $false, "415÷6=69, 1", 2) | Out-Null
$d.Content.Find.Execute("433÷3=144, 1", $true, $false, $false, $false, $false, $true, 1, $false, "171÷9=19, 0", 2) | Out-Null
$d.Content.Find.Execute("858÷6=143, 0", $true, $false, $false, $false, $false, $true, 1, $false, "339÷6=56, 3", 2) | Out-Null
$d.Content.Find.Execute("643÷8=80, 3", $true, $false, $false, $false, $false, $true, 1, $false, "794÷3=264, 2", 2) | Out-Null
$d.Content.Find.Execute("676÷2=338, 0", $true, $false, $false, $false, $false, $true, 1, $false, "540÷3=180, 0", 2) | Out-Null
$d.Content.Find.Execute("783÷6=130, 3", $true, $false, $false, $false, $false, $true, 1, $false, "717÷5=143, 2", 2) | Out-Null
$d.Content.Find.Execute("274÷6=45, 4", $true, $false, $false, $false, $false, $true, 1, $false, "133÷4=33, 1", 2) | Out-Null
$d.Content.Find.Execute("496÷9=55, 1", $true, $false, $false, $false, $false, $true, 1, $false, "852÷6=142, 0", 2) | Out-Null
$d.Content.Find.Execute("744÷2=372, 0", $true, $false, $false, $false, $false, $true, 1, $false, "702÷3=234, 0", 2) | Out-Null
$d.Content.Find.Execute("644÷2=322, 0", $true, $false, $false, $false, $false, $true, 1, $false, "204÷8=25, 4", 2) | Out-Null
$d.Content.Find.Execute("810÷6=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "955÷5=191, 0", 2) | Out-Null
$d.Content.Find.Execute("149÷7=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "837÷6=139, 3", 2) | Out-Null
$d.Content.Find.Execute("889÷9=98, 7", $true, $false, $false, $false, $false, $true, 1, $false, "576÷9=64, 0", 2) | Out-Null
$d.Content.Find.Execute("348÷4=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "732÷2=366, 0", 2) | Out-Null
$d.Content.Find.Execute("655÷6=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "826÷4=206, 2", 2) | Out-Null
$d.Content.Find.Execute("774÷4=193, 2", $true, $false, $false, $false, $false, $true, 1, $false, "436÷6=72, 4", 2) | Out-Null
$d.Content.Find.Execute("674÷9=74, 8", $true, $false, $false, $false, $false, $true, 1, $false, "988÷4=247, 0", 2) | Out-Null
$d.Content.Find.Execute("803÷4=200, 3", $true, $false, $false, $false, $false, $true, 1, $false, "634÷7=90, 4", 2) | Out-Null
$d.Content.Find.Execute("875÷9=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "603÷7=86, 1", 2) | Out-Null
$d.Content.Find.Execute("345÷8=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "238÷3=79, 1", 2) | Out-Null
$d.Content.Find.Execute("482÷8=60, 2", $true, $false, $false, $false, $false, $true, 1, $false, "160÷3=53, 1", 2) | Out-Null
